# Apply "Update GetVoteSummaries and SaveQSet (QSet table add column)" changes.
# This inserts 4 new rows into the QSets section of the error-code table
# (rows 149-152, pushing everything below down by 4), fills in the
# B/C (ErrorMessage/SP) columns for the existing QSets rows (1401-1406),
# and adds three brand-new QSets rows (1407, 1408, 1409).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Insert 4 rows right after row 148 (i.e. before the old row 149, which
# was blank separator before the "Qslides" header). This pushes the
# "Qslides" header (previously row 150) down to row 154, and everything
# below it shifts down by 4 rows to match the target layout.
$ws.Rows("149:152").Insert()

# --- Fill B/C columns for existing QSets rows (1401-1406) ---
$ws.Cells.Item(143, 2).Value = "Customer Id cannot be null or empty string."
$ws.Cells.Item(143, 3).Value = "SaveQSet"

$ws.Cells.Item(144, 2).Value = "Customer Id is not found."
$ws.Cells.Item(144, 3).Value = "SaveQSet"

$ws.Cells.Item(145, 2).Value = "QSet Id is not found."
$ws.Cells.Item(145, 3).Value = "SaveQSet"

$ws.Cells.Item(146, 2).Value = "QSet is already used in vote table."
$ws.Cells.Item(146, 3).Value = "SaveQSet"

$ws.Cells.Item(147, 2).Value = "Begin Date and/or End Date should not be null."
$ws.Cells.Item(147, 3).Value = "SaveQSet"

$ws.Cells.Item(148, 2).Value = "Display Mode is null or value is not in 0 to 1."
$ws.Cells.Item(148, 3).Value = "SaveQSet"

# --- New QSets rows (1407, 1408, 1409) ---
$ws.Cells.Item(149, 1).Value = 1407
$ws.Cells.Item(149, 2).Value = "Begin Date should less than End Date."
$ws.Cells.Item(149, 3).Value = "SaveQSet"

$ws.Cells.Item(150, 1).Value = 1408
$ws.Cells.Item(150, 2).Value = "Begin Date or End Date is overlap with another Question Set."
$ws.Cells.Item(150, 3).Value = "SaveQSet"

$ws.Cells.Item(151, 1).Value = 1409
$ws.Cells.Item(151, 3).Value = "SaveQSet"

# --- Fix up the sheet view to match the author's final cursor position ---
$ws.Range("B141").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 130
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
